$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rewrite the question rows (A2:A19) into clean numeric order with
#        normalized spacing. Assigning in increasing row order reproduces the
#        shared-string pool order of the target file. ---
$rowTexts = @{
  2  = "01. I wear blue-filtering,orange-tinted,and/or red-tinted glasses indoors during the day."
  3  = "02. I wear blue-filtering, orange-tinted, and/or red-tinted glasses outdoors during the day."
  4  = "03. I wear blue-filtering, orange-tinted, and/or red-tinted glasses within 1 hour before attempting to fall asleep."
  5  = "04. I spend 30 minutes or less per day (in total) outside."
  6  = "05. I spend between 30 minutes and 1 hour per day (in total) outside."
  7  = "06. I spend between 1 and 3 hours per day (in total) outside."
  8  = "07. I spend more than 3 hours per day (in total) outside."
  9  = "08. I spend as much time outside as possible."
  10 = "09. I go for a walk or exercise outside within 2 hours after waking up."
  11 = "10. I use my mobile phone within 1 hour before attempting to fall asleep."
  12 = "11. I look at my mobile phone screen immediately after waking up."
  13 = "12. I check my phone when I wake up at night."
  14 = "13. I dim my mobile phone screen within 1 hour before attempting to fall asleep."
  15 = "14. I use a blue-filter app on my computer screen within 1 hour before attempting to fall asleep."
  16 = "15. I dim my computer screen within 1 hour before attempting to fall asleep."
  17 = "16. I use tunable lights to create a healthy light environment."
  18 = "17. I use LEDs to create a healthy light environment."
  19 = "18. I use an alarm with a dawn simulation light."
}

for ($r = 2; $r -le 19; $r++) {
  $ws.Cells.Item($r, 1).Value = $rowTexts[$r]
}

# --- 2. Rows 2 & 3 no longer carry the big "section title" look (Arial 20) ---
#        they fall back to the workbook's default cell style.
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(4).RowHeight
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(4).RowHeight

# --- 3. Header cell A1 ("Items") becomes bold + centered ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter

# --- 4. Two new empty "section title" rows appended at 22 & 23, carrying the
#        big Arial-20 look that used to live on rows 2 & 3. ---
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Rows.Item(22).RowHeight = 25
$ws.Rows.Item(23).RowHeight = 25

# --- 5. View tweaks: zoom to 140%, move the active selection to A16 ---
$excel.ActiveWindow.Zoom = 140
$ws.Range("A16").Select()
